# Applies the "xong co ban ve dem va chinh sua fix 1 so thu" edit:
#  - tweak a couple of existing names
#  - fix a unit abbreviation
#  - add three new VDV rows (7,8,9) reusing the existing table formatting
#  - drop the now-unused trailing blank rows (17-19) so the table ends at row 16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix a couple of existing entries -----------------------------------
$ws.Range("B4").Value = "Thầy Toàn Cờ Vua 1"
$ws.Range("B5").Value = "Lê Trọng Đề Toàn 2"
$ws.Range("C8").Value = "hh"

# --- append new VDV rows, reusing row 9's look (font/border/alignment) --
$ws.Range("A9:C9").Copy()

$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Chim Sẻ Đi Nắng"
$ws.Range("C10").Value = "Hà Nội"
$ws.Range("A10:C10").PasteSpecial(-4122)

$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "Phùng Thị Tuyết Lan"
$ws.Range("C11").Value = "Hồ Chí Minh"
$ws.Range("A11:C11").PasteSpecial(-4122)

$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "Lê Quang Liêm"
$ws.Range("C12").Value = "HCM"
$ws.Range("A12:C12").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- drop the now superfluous trailing blank rows (19,18,17) ------------
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(18).Delete()
$ws.Rows.Item(17).Delete()

# --- restore the cursor/selection like the saved file shows -------------
$ws.Range("I16").Select()
